$p = $ppt.ActivePresentation

# Slide 17: 音声認識
$s17 = $p.Slides.Add(17, 2)
$s17.Shapes.Item(1).TextFrame.TextRange.Text = "音声認識"

$tr = $s17.Shapes.Item(2).TextFrame.TextRange
$tr.Text = "音声についても、特徴量が得られれば分類できる"
$tr = $tr.InsertAfter("`r使用例：")
$tr = $tr.InsertAfter("`rMFCC")
$tr = $tr.InsertAfter("（")
$tr = $tr.InsertAfter("Mel-Frequency Cepstral Coefficients")
$tr = $tr.InsertAfter("：メル周波数ケプストラム係数）")
$tr = $tr.InsertAfter("`r→音の「質感」や「母音らしさ」を数値化")
$tr = $tr.InsertAfter("`r`r")
$tr = $tr.InsertAfter("これをカーネルいれて、母音判定器の作成")

# Slide 18: 音声認識例
$s18 = $p.Slides.Add(18, 2)
$s18.Shapes.Item(1).TextFrame.TextRange.Text = "音声認識例"

$tf18 = $s18.Shapes.Item(2).TextFrame
$tf18.AutoSize = 2
$tr2 = $tf18.TextRange
$tr2.Text = "A,I,U,E,O"
$tr2 = $tr2.InsertAfter(" ")
$tr2 = $tr2.InsertAfter("の音声データを５つずつ訓練セットとして用意")
$tr2 = $tr2.InsertAfter("`r`rMFCC")
$tr2 = $tr2.InsertAfter("で特徴抽出")
$tr2 = $tr2.InsertAfter("`rLinear")
$tr2 = $tr2.InsertAfter("で分類")
$tr2 = $tr2.InsertAfter("`r結果")
$tr2 = $tr2.InsertAfter("`r🎤 テスト音声 テストい")
$tr2 = $tr2.InsertAfter(".wav ")
$tr2 = $tr2.InsertAfter("の予測結果")
$tr2 = $tr2.InsertAfter(": ")
$tr2 = $tr2.InsertAfter("i")
$tr2 = $tr2.InsertAfter("`r🎤 テスト音声 テストう")
$tr2 = $tr2.InsertAfter(".wav ")
$tr2 = $tr2.InsertAfter("の予測結果")
$tr2 = $tr2.InsertAfter(": u")
$tr2 = $tr2.InsertAfter("`r🎤 テスト音声 テストあ")
$tr2 = $tr2.InsertAfter(".wav ")
$tr2 = $tr2.InsertAfter("の予測結果")
$tr2 = $tr2.InsertAfter(": a")
$tr2 = $tr2.InsertAfter("`r🎤 テスト音声 テストえ")
$tr2 = $tr2.InsertAfter(".wav ")
$tr2 = $tr2.InsertAfter("の予測結果")
$tr2 = $tr2.InsertAfter(": e")
$tr2 = $tr2.InsertAfter("`r🎤 テスト音声 テストお")
$tr2 = $tr2.InsertAfter(".wav ")
$tr2 = $tr2.InsertAfter("の予測結果")
$tr2 = $tr2.InsertAfter(": o")

Write-Output ("Final slide count: " + $p.Slides.Count)
